# Generate Report for Archive
# Re-sorts the localization-status rows so that files currently "In
# Translation" (648008c6, 887e8d4e, b6321685) are grouped at the top
# (648008c6 moved from "Ready for handoff" -> "In Translation" and is
# now listed first), while 69f88810 / .localization-config stay put.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A=File Name, B=zh-cn status, C=de-de status
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Hyperlinks.Item(1).TextToDisplay = "648008c6-563f-4240-af01-96c9990c2be3.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "887e8d4e-e235-4a6c-b89d-db7b246b1693.md"
$ws.Hyperlinks.Item(3).TextToDisplay = "b6321685-eea9-4b0e-bdd3-a6985c0977a6.md"

$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Hyperlinks.Item(1).TextToDisplay = "648008c6-563f-4240-af01-96c9990c2be3.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "648008c6-563f-4240-af01-96c9990c2be3.937cedb4757a8a12eabb5ebe170d040ac2e5c2e9.zh-cn.xlf"
$ws.Hyperlinks.Item(3).TextToDisplay = "887e8d4e-e235-4a6c-b89d-db7b246b1693.md"
$ws.Hyperlinks.Item(4).TextToDisplay = "887e8d4e-e235-4a6c-b89d-db7b246b1693.ca66e44f0bf37f01ce4e46e8dec81c61ee965a34.zh-cn.xlf"
$ws.Hyperlinks.Item(5).TextToDisplay = "b6321685-eea9-4b0e-bdd3-a6985c0977a6.md"
$ws.Hyperlinks.Item(6).TextToDisplay = "b6321685-eea9-4b0e-bdd3-a6985c0977a6.2178b2e0e63e39bd299dcb9fa4bc4367aedc9be3.zh-cn.xlf"

$ws.Range("B2").Value = "In Translation"
$ws.Range("D2").Value = "2016-03-04 05:28:31"
$ws.Range("B3").Value = "In Translation"
$ws.Range("D3").Value = "2016-03-04 05:26:08"
$ws.Range("B4").Value = "In Translation"
$ws.Range("D4").Value = "2016-03-04 05:26:08"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Hyperlinks.Item(1).TextToDisplay = "648008c6-563f-4240-af01-96c9990c2be3.md"
$ws.Hyperlinks.Item(2).TextToDisplay = "648008c6-563f-4240-af01-96c9990c2be3.937cedb4757a8a12eabb5ebe170d040ac2e5c2e9.de-de.xlf"
$ws.Hyperlinks.Item(3).TextToDisplay = "887e8d4e-e235-4a6c-b89d-db7b246b1693.md"
$ws.Hyperlinks.Item(4).TextToDisplay = "887e8d4e-e235-4a6c-b89d-db7b246b1693.ca66e44f0bf37f01ce4e46e8dec81c61ee965a34.de-de.xlf"
$ws.Hyperlinks.Item(5).TextToDisplay = "b6321685-eea9-4b0e-bdd3-a6985c0977a6.md"
$ws.Hyperlinks.Item(6).TextToDisplay = "b6321685-eea9-4b0e-bdd3-a6985c0977a6.2178b2e0e63e39bd299dcb9fa4bc4367aedc9be3.de-de.xlf"

$ws.Range("B2").Value = "In Translation"
$ws.Range("D2").Value = "2016-03-04 05:28:47"
$ws.Range("B3").Value = "In Translation"
$ws.Range("D3").Value = "2016-03-04 05:26:40"
$ws.Range("B4").Value = "In Translation"
$ws.Range("D4").Value = "2016-03-04 05:26:40"
